$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '61.927.01'
$ws.Range("E2").Value = '  +1.57%  '

$ws.Range("D3").Value = '3.460.56'
$ws.Range("E3").Value = '  +3.58%  '

$ws.Range("E4").Value = '  -0.05%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '580.51'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +2.06%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '148.29'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +9.72%  '

$ws.Range("D7").Value = '3.460.90'
$ws.Range("E7").Value = '  +3.68%  '

$ws.Range("E8").Value = '  +0.04%  '

$ws.Range("E9").Value = '  +1.36%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '7.70'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +3.77%  '

$ws.Range("E11").Value = '  +2.63%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.388'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +0.88%  '

$ws.Range("D13").Value = '4.047.49'
$ws.Range("E13").Value = '  +3.53%  '

$ws.Range("B14").Value = 'TRON'
$ws.Range("C14").Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.123'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -0.32%  '

$ws.Range("B15").Value = 'Avalanche'
$ws.Range("C15").Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '27.90'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +8.40%  '

$ws.Range("E16").Value = '  +2.14%  '

$ws.Range("D17").Value = '3.457.74'
$ws.Range("E17").Value = '  +3.59%  '

$ws.Range("D18").Value = '61.940.69'
$ws.Range("E18").Value = '  +1.45%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '6.32'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +8.76%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '14.19'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +1.75%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '9.46'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +2.70%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '384.60'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +2.10%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.566'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +3.10%  '

$ws.Range("D24").Value = '3.587.74'
$ws.Range("E24").Value = '  +2.86%  '

$ws.Range("B25").Value = 'Dai'
$ws.Range("C25").Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '1.00'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +0.16%  '

$ws.Range("B26").Value = 'LEO'
$ws.Range("C26").Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '5.78'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +1.09%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '72.45'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +2.30%  '

$ws.Range("E28").Value = '  +1.21%  '

$ws.Range("E29").Value = '  +9.70%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '7.82'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +5.50%  '

$ws.Range("E31").Value = '  -10.81%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '1.00'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -0.12%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '8.23'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +1.54%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '2.18'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +2.56%  '

$ws.Range("E35").Value = '  -0.07%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '24.02'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +2.55%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '7.05'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +4.58%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '5.22'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +0.66%  '

$ws.Range("E39").Value = '  +3.02%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '166.76'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +1.50%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.0785'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +4.08%  '

$ws.Range("B42").Value = 'Mantle'
$ws.Range("C42").Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.798'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +4.27%  '

$ws.Range("B43").Value = 'EnergySwap'
$ws.Range("C43").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '25.91'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +9.24%  '

$ws.Range("E44").Value = '  +3.07%  '

$ws.Range("E45").Value = '  -0.04%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '42.33'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +2.51%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '4.49'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +2.83%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.18'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -1.55%  '

$ws.Range("D49").Value = '2.602.68'
$ws.Range("E49").Value = '  +11.33%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '23.57'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +3.59%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '6.87'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +1.39%  '
